# Update the Metadata worksheet and insert a "Jurisdiction" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "0.1.1"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row after row 10 (Contact) for Jurisdiction
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Copy style from the row above (Contact row, now row 10) to keep formatting consistent
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
